$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.79785166666667
$ws.Range("H2").Value = 47.39355500000001
$ws.Range("I2").Value = 0.1445757693628457
$ws.Range("J2").Value = 0.1445757693628457
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.776179
$ws.Range("N2").Value = 5.328537
$ws.Range("O2").Value = 0.009213114886297067
$ws.Range("P2").Value = 0.009213114886297067
$ws.Range("Q2").Value = 28.05981237544833
$ws.Range("R2").Value = 252.538311379035
$ws.Range("S2").Value = 0.001331993172914685
$ws.Range("T2").Value = 0.001331993172914685
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.79785166666667
$ws.Range("H3").Value = 47.39355500000001
$ws.Range("I3").Value = 0.1445757693628457
$ws.Range("J3").Value = 0.1445757693628457
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 103.273595
$ws.Range("N3").Value = 309.820785
$ws.Range("O3").Value = 0.5356844639284184
$ws.Range("P3").Value = 0.5356844639284185
$ws.Range("Q3").Value = 1631.500934893408
$ws.Range("R3").Value = 14683.50841404068
$ws.Range("S3").Value = 0.07744699350817466
$ws.Range("T3").Value = 0.0774469935081747
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.79785166666667
$ws.Range("H4").Value = 47.39355500000001
$ws.Range("I4").Value = 0.1445757693628457
$ws.Range("J4").Value = 0.1445757693628457
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 53.963124
$ws.Range("N4").Value = 161.889372
$ws.Range("O4").Value = 0.2799089849815219
$ws.Range("P4").Value = 0.2799089849815219
$ws.Range("Q4").Value = 852.5014284219401
$ws.Range("R4").Value = 7672.512855797461
$ws.Range("S4").Value = 0.04046805685527675
$ws.Range("T4").Value = 0.04046805685527676
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.79785166666667
$ws.Range("H5").Value = 47.39355500000001
$ws.Range("I5").Value = 0.1445757693628457
$ws.Range("J5").Value = 0.1445757693628457
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.77521133333333
$ws.Range("N5").Value = 101.325634
$ws.Range("O5").Value = 0.1751934362037625
$ws.Range("P5").Value = 0.1751934362037625
$ws.Range("Q5").Value = 533.5757786543189
$ws.Range("R5").Value = 4802.18200788887
$ws.Range("S5").Value = 0.02532872582647959
$ws.Range("T5").Value = 0.0253287258264796
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 74.39645633333333
$ws.Range("H6").Value = 223.189369
$ws.Range("I6").Value = 0.6808473164079603
$ws.Range("J6").Value = 0.6808473164079603
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.776179
$ws.Range("N6").Value = 5.328537
$ws.Range("O6").Value = 0.009213114886297067
$ws.Range("P6").Value = 0.009213114886297067
$ws.Range("Q6").Value = 132.1414234136837
$ws.Range("R6").Value = 1189.272810723153
$ws.Range("S6").Value = 0.006272724546093587
$ws.Range("T6").Value = 0.006272724546093587
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 74.39645633333333
$ws.Range("H7").Value = 223.189369
$ws.Range("I7").Value = 0.6808473164079603
$ws.Range("J7").Value = 0.6808473164079603
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 103.273595
$ws.Range("N7").Value = 309.820785
$ws.Range("O7").Value = 0.5356844639284184
$ws.Range("P7").Value = 0.5356844639284185
$ws.Range("Q7").Value = 7683.189500803852
$ws.Range("R7").Value = 69148.70550723467
$ws.Range("S7").Value = 0.3647193297071005
$ws.Range("T7").Value = 0.3647193297071005
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 74.39645633333333
$ws.Range("H8").Value = 223.189369
$ws.Range("I8").Value = 0.6808473164079603
$ws.Range("J8").Value = 0.6808473164079603
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 53.963124
$ws.Range("N8").Value = 161.889372
$ws.Range("O8").Value = 0.2799089849815219
$ws.Range("P8").Value = 0.2799089849815219
$ws.Range("Q8").Value = 4014.665198276252
$ws.Range("R8").Value = 36131.98678448627
$ws.Range("S8").Value = 0.1905752812631452
$ws.Range("T8").Value = 0.1905752812631452
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 74.39645633333333
$ws.Range("H9").Value = 223.189369
$ws.Range("I9").Value = 0.6808473164079603
$ws.Range("J9").Value = 0.6808473164079603
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.77521133333333
$ws.Range("N9").Value = 101.325634
$ws.Range("O9").Value = 0.1751934362037625
$ws.Range("P9").Value = 0.1751934362037625
$ws.Range("Q9").Value = 2512.756035109438
$ws.Range("R9").Value = 22614.80431598494
$ws.Range("S9").Value = 0.1192799808916209
$ws.Range("T9").Value = 0.1192799808916209
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 16.36992
$ws.Range("H10").Value = 49.10976
$ws.Range("I10").Value = 0.1498111153557632
$ws.Range("J10").Value = 0.1498111153557632
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.776179
$ws.Range("N10").Value = 5.328537
$ws.Range("O10").Value = 0.009213114886297067
$ws.Range("P10").Value = 0.009213114886297067
$ws.Range("Q10").Value = 29.07590813568
$ws.Range("R10").Value = 261.68317322112
$ws.Range("S10").Value = 0.001380227017016949
$ws.Range("T10").Value = 0.001380227017016949
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 16.36992
$ws.Range("H11").Value = 49.10976
$ws.Range("I11").Value = 0.1498111153557632
$ws.Range("J11").Value = 0.1498111153557632
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 103.273595
$ws.Range("N11").Value = 309.820785
$ws.Range("O11").Value = 0.5356844639284184
$ws.Range("P11").Value = 0.5356844639284185
$ws.Range("Q11").Value = 1690.5804882624
$ws.Range("R11").Value = 15215.2243943616
$ws.Range("S11").Value = 0.08025148701987043
$ws.Range("T11").Value = 0.08025148701987046
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 16.36992
$ws.Range("H12").Value = 49.10976
$ws.Range("I12").Value = 0.1498111153557632
$ws.Range("J12").Value = 0.1498111153557632
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 53.963124
$ws.Range("N12").Value = 161.889372
$ws.Range("O12").Value = 0.2799089849815219
$ws.Range("P12").Value = 0.2799089849815219
$ws.Range("Q12").Value = 883.37202283008
$ws.Range("R12").Value = 7950.348205470721
$ws.Range("S12").Value = 0.04193347723818135
$ws.Range("T12").Value = 0.04193347723818135
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 16.36992
$ws.Range("H13").Value = 49.10976
$ws.Range("I13").Value = 0.1498111153557632
$ws.Range("J13").Value = 0.1498111153557632
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.77521133333333
$ws.Range("N13").Value = 101.325634
$ws.Range("O13").Value = 0.1751934362037625
$ws.Range("P13").Value = 0.1751934362037625
$ws.Range("Q13").Value = 552.89750750976
$ws.Range("R13").Value = 4976.077567587839
$ws.Range("S13").Value = 0.0262459240806944
$ws.Range("T13").Value = 0.02624592408069441
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.706168666666667
$ws.Range("H14").Value = 8.118506
$ws.Range("I14").Value = 0.02476579887343077
$ws.Range("J14").Value = 0.02476579887343077
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.776179
$ws.Range("N14").Value = 5.328537
$ws.Range("O14").Value = 0.009213114886297067
$ws.Range("P14").Value = 0.009213114886297067
$ws.Range("Q14").Value = 4.806639956191333
$ws.Range("R14").Value = 43.259759605722
$ws.Range("S14").Value = 0.0002281701502718441
$ws.Range("T14").Value = 0.0002281701502718441
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.706168666666667
$ws.Range("H15").Value = 8.118506
$ws.Range("I15").Value = 0.02476579887343077
$ws.Range("J15").Value = 0.02476579887343077
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 103.273595
$ws.Range("N15").Value = 309.820785
$ws.Range("O15").Value = 0.5356844639284184
$ws.Range("P15").Value = 0.5356844639284185
$ws.Range("Q15").Value = 279.4757668830233
$ws.Range("R15").Value = 2515.28190194721
$ws.Range("S15").Value = 0.01326665369327279
$ws.Range("T15").Value = 0.01326665369327279
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.706168666666667
$ws.Range("H16").Value = 8.118506
$ws.Range("I16").Value = 0.02476579887343077
$ws.Range("J16").Value = 0.02476579887343077
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 53.963124
$ws.Range("N16").Value = 161.889372
$ws.Range("O16").Value = 0.2799089849815219
$ws.Range("P16").Value = 0.2799089849815219
$ws.Range("Q16").Value = 146.033315324248
$ws.Range("R16").Value = 1314.299837918232
$ws.Range("S16").Value = 0.006932169624918524
$ws.Range("T16").Value = 0.006932169624918524
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.706168666666667
$ws.Range("H17").Value = 8.118506
$ws.Range("I17").Value = 0.02476579887343077
$ws.Range("J17").Value = 0.02476579887343077
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.77521133333333
$ws.Range("N17").Value = 101.325634
$ws.Range("O17").Value = 0.1751934362037625
$ws.Range("P17").Value = 0.1751934362037625
$ws.Range("Q17").Value = 91.40141862031155
$ws.Range("R17").Value = 822.6127675828039
$ws.Range("S17").Value = 0.004338805404967607
$ws.Range("T17").Value = 0.004338805404967607